$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing used range to remove stale row 11 / stale columns beyond M
$ws.Cells.Clear()

# Un-ignore the "number stored as text" warning on the old range
$oldIgnored = $ws.Range("A1:H11").Errors.Item(9)
$oldIgnored.Ignore = $false

# Header row
$ws.Range("A1").Value = 'id'
$ws.Range("B1").Value = 'type'
$ws.Range("C1").Value = 'dc\.identifier'
$ws.Range("D1").Value = 'inScheme'
$ws.Range("E1").Value = 'member'
$ws.Range("F1").Value = 'prefLabel'
$ws.Range("G1").Value = 'seeAlso'
$ws.Range("H1").Value = 'definition'
$ws.Range("I1").Value = 'notation'
$ws.Range("J1").Value = 'note'
$ws.Range("K1").Value = 'topConceptOf'
$ws.Range("L1").Value = 'altLabel'
$ws.Range("M1").Value = 'hasTopConcept'

# Row 2
$ws.Range("A2").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/collection/consistentietype/consistentietypes'
$ws.Range("B2").Value = 'http://www.w3.org/2004/02/skos/core#Collection'
$ws.Range("C2").Value = 'be.vlaanderen.bodemenondergrond.data.id.collection.consistentietype.consistentietypes'
$ws.Range("D2").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("E2").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/plastisch|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/slap|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/stevig|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/stijf|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/versteend|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/zeerslap|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/zeerstijf'
$ws.Range("F2").Value = 'Collectie van consistentietypes.'
$ws.Range("G2").Value = 'null'
$ws.Range("H2").Value = 'null'
$ws.Range("I2").Value = 'null'
$ws.Range("J2").Value = 'null'
$ws.Range("K2").Value = 'null'
$ws.Range("L2").Value = 'null'
$ws.Range("M2").Value = 'null'

# Row 3
$ws.Range("A3").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/plastisch'
$ws.Range("B3").Value = 'http://www.w3.org/2004/02/skos/core#Concept'
$ws.Range("C3").Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.consistentietype.plastisch'
$ws.Range("D3").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("E3").Value = 'null'
$ws.Range("F3").Value = 'plastisch'
$ws.Range("G3").Value = 'https://app.nbn.be/data/r/platform/frontend/detail?p40_id=189507&p40_language_code=nl&p40_detail_id=85507'
$ws.Range("H3").Value = 'De grond laat toe om eenvoudig een rolletje van 3 mm  te maken.'
$ws.Range("I3").Value = 'plastisch'
$ws.Range("J3").Value = 'De grond laat toe om eenvoudig een rolletje van 3 mm  te maken.|In de norm NBN EN ISO 14688-1:2018 wordt er een onderscheid gemaakt tussen de eigenschappen consistentie en plasticiteit, met elk een eigen codelijst. Dit veld is behouden in deze versie van codetabel ''Consistentie'', maar moet gemigreerd worden naar een aparte tabel ''Plasticiteit'' conform de vermelde norm.'
$ws.Range("K3").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("L3").Value = 'null'
$ws.Range("M3").Value = 'null'

# Row 4
$ws.Range("A4").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/slap'
$ws.Range("B4").Value = 'http://www.w3.org/2004/02/skos/core#Concept'
$ws.Range("C4").Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.consistentietype.slap'
$ws.Range("D4").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("E4").Value = 'null'
$ws.Range("F4").Value = 'slap'
$ws.Range("G4").Value = 'https://app.nbn.be/data/r/platform/frontend/detail?p40_id=189507&p40_language_code=nl&p40_detail_id=85507'
$ws.Range("H4").Value = 'De grond kan tot 10 mm diep met de vinger worden ingedrukt en kan met lichte druk van de vingers worden verkneed.'
$ws.Range("I4").Value = 'slap'
$ws.Range("J4").Value = 'De grond kan tot 10 mm diep met de vinger worden ingedrukt en kan met lichte druk van de vingers worden verkneed.|Indeling afgeleid van NBN EN ISO 14688-1:2018 en vooral van toepassing voor laag (weinig tot matig) plastische gronden.'
$ws.Range("K4").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("L4").Value = 'null'
$ws.Range("M4").Value = 'null'

# Row 5
$ws.Range("A5").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/stevig'
$ws.Range("B5").Value = 'http://www.w3.org/2004/02/skos/core#Concept'
$ws.Range("C5").Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.consistentietype.stevig'
$ws.Range("D5").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("E5").Value = 'null'
$ws.Range("F5").Value = 'stevig'
$ws.Range("G5").Value = 'https://app.nbn.be/data/r/platform/frontend/detail?p40_id=189507&p40_language_code=nl&p40_detail_id=85507'
$ws.Range("H5").Value = 'null'
$ws.Range("I5").Value = 'stevig'
$ws.Range("J5").Value = 'De grond kan met de hand gerold worden tot 3 mm dikke strengen zonder te breken of te verkruimelen. Het kan gemakkelijk met de duim worden ingedrukt, maar kan niet met de vingers worden verkneed.'
$ws.Range("K5").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("L5").Value = 'matig vast'
$ws.Range("M5").Value = 'null'

# Row 6
$ws.Range("A6").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/stijf'
$ws.Range("B6").Value = 'http://www.w3.org/2004/02/skos/core#Concept'
$ws.Range("C6").Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.consistentietype.stijf'
$ws.Range("D6").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("E6").Value = 'null'
$ws.Range("F6").Value = 'stijf'
$ws.Range("G6").Value = 'https://app.nbn.be/data/r/platform/frontend/detail?p40_id=189507&p40_language_code=nl&p40_detail_id=85507'
$ws.Range("H6").Value = 'null'
$ws.Range("I6").Value = 'stijf'
$ws.Range("J6").Value = 'De grond kan niet met de vingers worden verkneed en verkruimelt of breekt wanneer deze tot 3 mm dikke strengen wordt gerold, maar is nog vochtig genoeg om weer tot een bol te worden gekneed. Met de duim kan een ondiepe indrukking gemaakt worden.'
$ws.Range("K6").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("L6").Value = 'vast'
$ws.Range("M6").Value = 'null'

# Row 7
$ws.Range("A7").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/versteend'
$ws.Range("B7").Value = 'http://www.w3.org/2004/02/skos/core#Concept'
$ws.Range("C7").Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.consistentietype.versteend'
$ws.Range("D7").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("E7").Value = 'null'
$ws.Range("F7").Value = 'versteend'
$ws.Range("G7").Value = 'https://app.nbn.be/data/r/platform/frontend/detail?p40_id=189507&p40_language_code=nl&p40_detail_id=85507'
$ws.Range("H7").Value = 'De grond is versteend en kan niet met de hand in water worden gedesaggregeerd. Het heeft een natuurlijke samenstelling of bestaat uit een aggregaat van mineraalkorrels, kristallen of op mineralen gebaseerde deeltjes, die zijn samengeperst, gecementeerd of anderszins aan elkaar zijn gebonden.'
$ws.Range("I7").Value = 'versteend'
$ws.Range("J7").Value = 'Indeling afgeleid van NBN EN ISO 14688-1:2018 en vooral van toepassing voor laag (weinig tot matig) plastische gronden.|De grond is versteend en kan niet met de hand in water worden gedesaggregeerd. Het heeft een natuurlijke samenstelling of bestaat uit een aggregaat van mineraalkorrels, kristallen of op mineralen gebaseerde deeltjes, die zijn samengeperst, gecementeerd of anderszins aan elkaar zijn gebonden.'
$ws.Range("K7").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("L7").Value = 'null'
$ws.Range("M7").Value = 'null'

# Row 8
$ws.Range("A8").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/zeerslap'
$ws.Range("B8").Value = 'http://www.w3.org/2004/02/skos/core#Concept'
$ws.Range("C8").Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.consistentietype.zeerslap'
$ws.Range("D8").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("E8").Value = 'null'
$ws.Range("F8").Value = 'zeer slap'
$ws.Range("G8").Value = 'https://app.nbn.be/data/r/platform/frontend/detail?p40_id=189507&p40_language_code=nl&p40_detail_id=85507'
$ws.Range("H8").Value = 'De grond kan tot 25 mm diep met de vinger worden ingedrukt en loopt tussen de vingers door wanneer het in de hand wordt samengeknepen.'
$ws.Range("I8").Value = 'zeerslap'
$ws.Range("J8").Value = 'Indeling afgeleid van NBN EN ISO 14688-1:2018 en vooral van toepassing voor laag (weinig tot matig) plastische gronden.|De grond kan tot 25 mm diep met de vinger worden ingedrukt en loopt tussen de vingers door wanneer het in de hand wordt samengeknepen.'
$ws.Range("K8").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("L8").Value = 'null'
$ws.Range("M8").Value = 'null'

# Row 9
$ws.Range("A9").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/zeerstijf'
$ws.Range("B9").Value = 'http://www.w3.org/2004/02/skos/core#Concept'
$ws.Range("C9").Value = 'be.vlaanderen.bodemenondergrond.data.id.concept.consistentietype.zeerstijf'
$ws.Range("D9").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("E9").Value = 'null'
$ws.Range("F9").Value = 'zeer stijf'
$ws.Range("G9").Value = 'https://app.nbn.be/data/r/platform/frontend/detail?p40_id=189507&p40_language_code=nl&p40_detail_id=85507'
$ws.Range("H9").Value = 'null'
$ws.Range("I9").Value = 'zeerstijf'
$ws.Range("J9").Value = 'De grond kan niet meer worden vervormd en verkruimelt onder druk.  Met de duimnagel kan de grond ingekerfd worden.'
$ws.Range("K9").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("L9").Value = 'hard'
$ws.Range("M9").Value = 'null'

# Row 10
$ws.Range("A10").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/consistentietype'
$ws.Range("B10").Value = 'http://www.w3.org/2004/02/skos/core#ConceptScheme'
$ws.Range("C10").Value = 'be.vlaanderen.bodemenondergrond.data.id.conceptscheme.consistentietype'
$ws.Range("D10").Value = 'null'
$ws.Range("E10").Value = 'null'
$ws.Range("F10").Value = 'Conceptschema consistentietype'
$ws.Range("G10").Value = 'null'
$ws.Range("H10").Value = 'null'
$ws.Range("I10").Value = 'null'
$ws.Range("J10").Value = 'null'
$ws.Range("K10").Value = 'null'
$ws.Range("L10").Value = 'null'
$ws.Range("M10").Value = 'https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/plastisch|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/slap|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/stevig|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/stijf|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/versteend|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/zeerslap|https://data.bodemenondergrond.vlaanderen.be/id/concept/consistentietype/zeerstijf'

# Re-apply the "number stored as text" ignored-error marker over the new used range
$newIgnored = $ws.Range("A1:M10").Errors.Item(9)
$newIgnored.Ignore = $true
